$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "calibration"
$ws.Range("G3").Value = "calibration"
$ws.Range("G4").Value = "deriv_calibration"
$ws.Range("G5").Value = "deriv_calibration"
$ws.Range("G6").Value = "calibration"
$ws.Range("G7").Value = "calibration"
